# Insert a new data row at row 715 (pushes 2026/12/29.. block, and everything
# after it, down by one row) and populate it with the new reading that the
# diff introduces: 2026/01/26, 月, 5, 201.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(715).Insert()

# Format column A as Text first so the "yyyy/mm/dd"-shaped string is stored
# as a literal string (matching the rest of the column) instead of being
# auto-converted into a date serial number by Excel's smart input parsing.
$ws.Range("A715").NumberFormat = "@"
$ws.Range("A715").Value = "2026/01/26"
$ws.Range("A715").ClearFormats()

$ws.Range("B715").Value = "月"
$ws.Range("C715").Value = 5
$ws.Range("D715").Value = 201
